$wb = $excel.ActiveWorkbook

# --- Staff: add new mentor/staff member "Jackie Ko" ---
$wsStaff = $wb.Worksheets.Item("Staff")
$wsStaff.Range("A3").Value = "Jackie Ko"
$wsStaff.Range("A3").Select()

# --- Mentor: add "Jackie Ko" from "SFSU" ---
$wsMentor = $wb.Worksheets.Item("Mentor")
$wsMentor.Range("A3").Value = "Jackie Ko"
$wsMentor.Range("B3").Value = "SFSU"
$wsMentor.Range("A3:B3").Select()

# --- Student: add "Jackie Ko" / "SFSU" row ---
$wsStudent = $wb.Worksheets.Item("Student")
$wsStudent.Range("A5").Value = "Jackie Ko"
$wsStudent.Range("B5").Value = "SFSU"
$wsStudent.Range("A5:B5").Select()

# --- Tool: add "Dewalt Drill" with quantity 7 ---
$wsTool = $wb.Worksheets.Item("Tool")
$wsTool.Range("A5").Value = "Dewalt Drill"
$wsTool.Range("B5").Value = 7
$wsTool.Range("A5:B5").Select()

# --- Purchase: add "Plus Gusset" purchase row, leave this the active sheet/tab ---
$wsPurchase = $wb.Worksheets.Item("Purchase")
$wsPurchase.Range("A5").Value = "Plus Gusset"
$wsPurchase.Range("B5").Value = 0.2
$wsPurchase.Range("C5").Value = 1000
$wsPurchase.Range("A5:C5").Select()
$wsPurchase.Activate()
